# Apply the recorded edit:
#  1. Fix the mis-typed header text "MODEL_CONDITION" -> "MODELCONDITION".
#  2. Delete column A entirely (the stray 9/12 values that were wrongly
#     carrying the header style). This shifts every other column one to
#     the left (B->A, C->B, D->C, E->D, F->E) and drops that leftover
#     style from the shifted-in data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# In-place text fix for the header label (keeps it the same cell/string
# rather than writing a brand-new value).
[void]$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# Deleting column A shifts B:F left to A:E.
$ws.Range("A:A").Delete()
